# Append new "completed breakout" rows to the three data sheets, as produced
# by the stock.yaml breakout job run on 10/06/2024.

$wb = $excel.ActiveWorkbook
$dateFmt = "YYYY-MM-DD HH:MM:SS"

# ---------------------------------------------------------------------------
# Sheet "three_line" -> new row 117
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("three_line")
$r = 117
$ws1.Cells.Item($r, 1).Value = 45447.55208333334
$ws1.Cells.Item($r, 1).NumberFormat = $dateFmt
$ws1.Cells.Item($r, 2).Value = "10-06-2024 10:15:00"
$ws1.Cells.Item($r, 3).Value = "hour"
$ws1.Cells.Item($r, 4).Value = "DWARKESH.NS"
$ws1.Cells.Item($r, 5).Value = 45418.38541666666
$ws1.Cells.Item($r, 5).NumberFormat = $dateFmt
$ws1.Cells.Item($r, 6).Value = 72.75
$ws1.Cells.Item($r, 7).Value = 45436.55208333334
$ws1.Cells.Item($r, 7).NumberFormat = $dateFmt
$ws1.Cells.Item($r, 8).Value = 71.90000152587891
$ws1.Cells.Item($r, 9).Value = 45446.38541666666
$ws1.Cells.Item($r, 9).NumberFormat = $dateFmt
$ws1.Cells.Item($r, 10).Value = 72
$ws1.Cells.Item($r, 11).Value = "High"
$ws1.Cells.Item($r, 12).Value = "10/06/2024 04:47:27"

# ---------------------------------------------------------------------------
# Sheet "two_line" -> new rows 14, 15
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("two_line")

$r = 14
$ws2.Cells.Item($r, 1).Value = 45435.38541666666
$ws2.Cells.Item($r, 1).NumberFormat = $dateFmt
$ws2.Cells.Item($r, 2).Value = "10-06-2024 10:15:00"
$ws2.Cells.Item($r, 3).Value = "hour"
$ws2.Cells.Item($r, 4).Value = "ACC.NS"
$ws2.Cells.Item($r, 5).Value = 45415.38541666666
$ws2.Cells.Item($r, 5).NumberFormat = $dateFmt
$ws2.Cells.Item($r, 6).Value = 2564
$ws2.Cells.Item($r, 7).Value = 45433.51041666666
$ws2.Cells.Item($r, 7).NumberFormat = $dateFmt
$ws2.Cells.Item($r, 8).Value = 2559.89990234375
$ws2.Cells.Item($r, 9).Value = "High"
$ws2.Cells.Item($r, 10).Value = "10/06/2024 04:47:27"

$r = 15
$ws2.Cells.Item($r, 1).Value = 45435.55208333334
$ws2.Cells.Item($r, 1).NumberFormat = $dateFmt
$ws2.Cells.Item($r, 2).Value = "10-06-2024 10:15:00"
$ws2.Cells.Item($r, 3).Value = "hour"
$ws2.Cells.Item($r, 4).Value = "BORORENEW.NS"
$ws2.Cells.Item($r, 5).Value = 45433.38541666666
$ws2.Cells.Item($r, 5).NumberFormat = $dateFmt
$ws2.Cells.Item($r, 6).Value = 495.6000061035156
$ws2.Cells.Item($r, 7).Value = 45434.38541666666
$ws2.Cells.Item($r, 7).NumberFormat = $dateFmt
$ws2.Cells.Item($r, 8).Value = 495.6000061035156
$ws2.Cells.Item($r, 9).Value = "Low"
$ws2.Cells.Item($r, 10).Value = "10/06/2024 04:47:27"

# ---------------------------------------------------------------------------
# Sheet "ph_pl_breakout_line" -> new rows 492-497
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("ph_pl_breakout_line")

$r = 492
$ws3.Cells.Item($r, 1).Value = "PICCADIL.BO"
$ws3.Cells.Item($r, 2).Value = 45446.38541666666
$ws3.Cells.Item($r, 2).NumberFormat = $dateFmt
$ws3.Cells.Item($r, 3).Value = 813.2999877929688
$ws3.Cells.Item($r, 4).Value = 750
$ws3.Cells.Item($r, 5).Value = 790
$ws3.Cells.Item($r, 6).Value = "High"
$ws3.Cells.Item($r, 7).Value = 813.2999877929688
$ws3.Cells.Item($r, 8).Value = "hour"
$ws3.Cells.Item($r, 9).Value = "10-06-2024 09:15:00"
$ws3.Cells.Item($r, 10).Value = 819
$ws3.Cells.Item($r, 11).Value = 800
$ws3.Cells.Item($r, 12).Value = "10/06/2024 04:47:27"

$r = 493
$ws3.Cells.Item($r, 1).Value = "TRIL.BO"
$ws3.Cells.Item($r, 2).Value = 45449.38541666666
$ws3.Cells.Item($r, 2).NumberFormat = $dateFmt
$ws3.Cells.Item($r, 3).Value = 750
$ws3.Cells.Item($r, 4).Value = 700
$ws3.Cells.Item($r, 5).Value = 722
$ws3.Cells.Item($r, 6).Value = "High"
$ws3.Cells.Item($r, 7).Value = 750
$ws3.Cells.Item($r, 8).Value = "hour"
$ws3.Cells.Item($r, 9).Value = "10-06-2024 09:15:00"
$ws3.Cells.Item($r, 10).Value = 764.4000244140625
$ws3.Cells.Item($r, 11).Value = 749.4000244140625
$ws3.Cells.Item($r, 12).Value = "10/06/2024 04:47:27"

$r = 494
$ws3.Cells.Item($r, 1).Value = "TRIL.BO"
$ws3.Cells.Item($r, 2).Value = 45449.42708333334
$ws3.Cells.Item($r, 2).NumberFormat = $dateFmt
$ws3.Cells.Item($r, 3).Value = 750
$ws3.Cells.Item($r, 4).Value = 728
$ws3.Cells.Item($r, 5).Value = 740
$ws3.Cells.Item($r, 6).Value = "High"
$ws3.Cells.Item($r, 7).Value = 750
$ws3.Cells.Item($r, 8).Value = "hour"
$ws3.Cells.Item($r, 9).Value = "10-06-2024 09:15:00"
$ws3.Cells.Item($r, 10).Value = 764.4000244140625
$ws3.Cells.Item($r, 11).Value = 749.4000244140625
$ws3.Cells.Item($r, 12).Value = "10/06/2024 04:47:27"

$r = 495
$ws3.Cells.Item($r, 1).Value = "JWL.NS"
$ws3.Cells.Item($r, 2).Value = 45449.38541666666
$ws3.Cells.Item($r, 2).NumberFormat = $dateFmt
$ws3.Cells.Item($r, 3).Value = 613
$ws3.Cells.Item($r, 4).Value = 563.3499755859375
$ws3.Cells.Item($r, 5).Value = 605.0499877929688
$ws3.Cells.Item($r, 6).Value = "High"
$ws3.Cells.Item($r, 7).Value = 613
$ws3.Cells.Item($r, 8).Value = "hour"
$ws3.Cells.Item($r, 9).Value = "10-06-2024 10:15:00"
$ws3.Cells.Item($r, 10).Value = 614
$ws3.Cells.Item($r, 11).Value = 612.4000244140625
$ws3.Cells.Item($r, 12).Value = "10/06/2024 04:47:27"

$r = 496
$ws3.Cells.Item($r, 1).Value = "KPEL.BO"
$ws3.Cells.Item($r, 2).Value = 45442.38541666666
$ws3.Cells.Item($r, 2).NumberFormat = $dateFmt
$ws3.Cells.Item($r, 3).Value = 437.5
$ws3.Cells.Item($r, 4).Value = 420.2000122070312
$ws3.Cells.Item($r, 5).Value = 430.25
$ws3.Cells.Item($r, 6).Value = "High"
$ws3.Cells.Item($r, 7).Value = 437.5
$ws3.Cells.Item($r, 8).Value = "hour"
$ws3.Cells.Item($r, 9).Value = "10-06-2024 09:15:00"
$ws3.Cells.Item($r, 10).Value = 438
$ws3.Cells.Item($r, 11).Value = 430.9500122070312
$ws3.Cells.Item($r, 12).Value = "10/06/2024 04:47:27"

$r = 497
$ws3.Cells.Item($r, 1).Value = "KPEL.BO"
$ws3.Cells.Item($r, 2).Value = 45446.38541666666
$ws3.Cells.Item($r, 2).NumberFormat = $dateFmt
$ws3.Cells.Item($r, 3).Value = 436.2000122070312
$ws3.Cells.Item($r, 4).Value = 423.7999877929688
$ws3.Cells.Item($r, 5).Value = 428.8999938964844
$ws3.Cells.Item($r, 6).Value = "High"
$ws3.Cells.Item($r, 7).Value = 436.2000122070312
$ws3.Cells.Item($r, 8).Value = "hour"
$ws3.Cells.Item($r, 9).Value = "10-06-2024 09:15:00"
$ws3.Cells.Item($r, 10).Value = 438
$ws3.Cells.Item($r, 11).Value = 430.9500122070312
$ws3.Cells.Item($r, 12).Value = "10/06/2024 04:47:27"
